$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1, G1)
$ws.Range("F1").Value = "pane number in tmux"
$ws.Range("G1").Value = "status"

# New "status" column values for existing rows 3 and 4
$ws.Range("G3").Value = "Done!"
$ws.Range("G4").Value = "Done!"

# New row 5 - TPR only experiment
$ws.Range("A5").Value = "Just TPR no LSTM in `nphrase embedding layer"
$ws.Range("B5").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR False --justTPR True --batch_size 60 --run_id 0 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP3.txt"
$ws.Range("C5").Value = "DLT1 / 3"
$ws.Range("D5").Value = "EXP3.txt"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("A5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 30

# New row 6 - LSTM only experiment
$ws.Range("A6").Value = "Just LSTM no TPR in `nphrase embedding layer"
$ws.Range("B6").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR False --justLSTM True --batch_size 60 --run_id 1 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP4.txt"
$ws.Range("C6").Value = "DLT1 / 4"
$ws.Range("D6").Value = "EXP4.txt"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("A6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

# Column F width
$ws.Columns.Item(6).ColumnWidth = 19.6

# Update selection to match target
$ws.Range("D10").Select()
